$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.74
$ws.Range("J2").Value = 3.8
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 5.1
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.46
$ws.Range("S2").Value = 2.64
$ws.Range("X2").Value = 24
$ws.Range("AH2").Value = 18.5
$ws.Range("AK2").Value = 32
$ws.Range("F3").Value = 2.12
$ws.Range("G3").Value = 2.36
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 3.25
$ws.Range("J3").Value = 3.75
$ws.Range("K3").Value = 5.9
$ws.Range("L3").Value = 1.01
$ws.Range("N3").Value = 3.45
$ws.Range("O3").Value = 1.12
$ws.Range("P3").Value = 2.78
$ws.Range("Q3").Value = 1.38
$ws.Range("R3").Value = 1.85
$ws.Range("S3").Value = 1.8
$ws.Range("T3").Value = 1.4
$ws.Range("U3").Value = 2.96
$ws.Range("V3").Value = 1.45
$ws.Range("W3").Value = 1.73
$ws.Range("AN3").Value = 10
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.34
$ws.Range("T4").Value = 2.06
$ws.Range("U4").Value = 1.79
$ws.Range("F6").Value = 3.45
$ws.Range("G6").Value = 3.85
$ws.Range("H6").Value = 2.14
$ws.Range("I6").Value = 2.32
$ws.Range("J6").Value = 3.2
$ws.Range("L6").Value = 1.34
$ws.Range("N6").Value = 3.6
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.9
$ws.Range("Q6").Value = 1.94
$ws.Range("R6").Value = 1.35
$ws.Range("S6").Value = 3.4
$ws.Range("T6").Value = 1.76
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.75
$ws.Range("W6").Value = 1.35
$ws.Range("X6").Value = 17.5
$ws.Range("Y6").Value = 12
$ws.Range("Z6").Value = 17.5
$ws.Range("AA6").Value = 36
$ws.Range("AB6").Value = 16.5
$ws.Range("AD6").Value = 13.5
$ws.Range("AE6").Value = 29
$ws.Range("AF6").Value = 32
$ws.Range("AG6").Value = 18.5
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 46
$ws.Range("AJ6").Value = 85
$ws.Range("AK6").Value = 55
$ws.Range("AL6").Value = 65
$ws.Range("AM6").Value = 120
$ws.Range("AO6").Value = 22
$ws.Range("F7").Value = 1.31
$ws.Range("J7").Value = 5.3
$ws.Range("P7").Value = 2.2
$ws.Range("Q7").Value = 1.48
$ws.Range("R7").Value = 1.56
$ws.Range("S7").Value = 2.22
$ws.Range("T7").Value = 1.94
$ws.Range("U7").Value = 1.84
$ws.Range("X7").Value = 30
$ws.Range("Y7").Value = 46
$ws.Range("AC7").Value = 17
$ws.Range("AD7").Value = 48
$ws.Range("AJ7").Value = 13.5
$ws.Range("AK7").Value = 17
$ws.Range("AL7").Value = 44
$ws.Range("AN7").Value = 5.8
$ws.Range("F8").Value = 2.78
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 2.5
$ws.Range("N8").Value = 4
$ws.Range("O8").Value = 1.27
$ws.Range("P8").Value = 2.04
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.4
$ws.Range("S8").Value = 3
$ws.Range("T8").Value = 1.67
$ws.Range("U8").Value = 2.22
$ws.Range("V8").Value = 1.59
$ws.Range("Y8").Value = 12.5
$ws.Range("AB8").Value = 13
$ws.Range("AC8").Value = 8.2
$ws.Range("AI8").Value = 38
$ws.Range("AJ8").Value = 55
$ws.Range("AL8").Value = 40
$ws.Range("AO8").Value = 21
$ws.Range("F9").Value = 3.8
$ws.Range("G9").Value = 4.2
$ws.Range("H9").Value = 1.89
$ws.Range("I9").Value = 2.02
$ws.Range("L9").Value = 1.25
$ws.Range("N9").Value = 5.4
$ws.Range("P9").Value = 2.5
$ws.Range("Q9").Value = 1.57
$ws.Range("R9").Value = 1.61
$ws.Range("S9").Value = 2.4
$ws.Range("T9").Value = 1.57
$ws.Range("U9").Value = 2.5
$ws.Range("V9").Value = 1.98
$ws.Range("W9").Value = 1.32
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 26
$ws.Range("AC9").Value = 10.5
$ws.Range("AD9").Value = 11
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 16
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 48
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 9
$ws.Range("F11").Value = 3.05
$ws.Range("G11").Value = 3.5
$ws.Range("H11").Value = 2.12
$ws.Range("I11").Value = 2.32
$ws.Range("J11").Value = 3.9
$ws.Range("K11").Value = 4.5
$ws.Range("T11").Value = 1.52
$ws.Range("V11").Value = 1.75
$ws.Range("W11").Value = 1.4
$ws.Range("I13").Value = 2.88
$ws.Range("U13").Value = 2.2
$ws.Range("AB13").Value = 11.5
$ws.Range("N14").Value = 3.15
$ws.Range("O14").Value = 1.39
$ws.Range("U14").Value = 1.58
$ws.Range("F15").Value = 3.6
$ws.Range("G15").Value = 3.95
$ws.Range("H15").Value = 2.38
$ws.Range("I15").Value = 2.54
$ws.Range("J15").Value = 2.94
$ws.Range("V15").Value = 1.64
$ws.Range("W15").Value = 1.34
$ws.Range("Y15").Value = 6.2
$ws.Range("Z15").Value = 13
$ws.Range("AA15").Value = 80
$ws.Range("AB15").Value = 8.4
$ws.Range("AC15").Value = 7.8
$ws.Range("AD15").Value = 17.5
$ws.Range("AE15").Value = 140
$ws.Range("AF15").Value = 24
$ws.Range("AG15").Value = 26
$ws.Range("AI15").Value = 120
$ws.Range("AL15").Value = 150
$ws.Range("AN15").Value = 170
$ws.Range("AO15").Value = 300
$ws.Range("G16").Value = 2.16
$ws.Range("H16").Value = 4.5
$ws.Range("I16").Value = 5.2
$ws.Range("J16").Value = 3.15
$ws.Range("K16").Value = 3.4
$ws.Range("O16").Value = 1.55
$ws.Range("S16").Value = 5.3
$ws.Range("T16").Value = 2.2
$ws.Range("U16").Value = 1.59
$ws.Range("V16").Value = 1.24
$ws.Range("W16").Value = 1.86
$ws.Range("Y16").Value = 1000
$ws.Range("AB16").Value = 1000
$ws.Range("O17").Value = 1.56
$ws.Range("AD17").Value = 12
$ws.Range("G18").Value = 1.86
$ws.Range("H18").Value = 4.4
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 3.6
$ws.Range("K18").Value = 4.6
$ws.Range("L18").Value = 1.43
$ws.Range("N18").Value = 3.25
$ws.Range("O18").Value = 1.39
$ws.Range("P18").Value = 1.76
$ws.Range("Q18").Value = 2
$ws.Range("S18").Value = 3.65
$ws.Range("T18").Value = 1.72
$ws.Range("W18").Value = 2.16
$ws.Range("F19").Value = 3.55
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 2.2
$ws.Range("I19").Value = 2.38
$ws.Range("J19").Value = 3.15
$ws.Range("K19").Value = 3.75
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 2.98
$ws.Range("P19").Value = 1.67
$ws.Range("Q19").Value = 2.26
$ws.Range("S19").Value = 3.8
$ws.Range("T19").Value = 1.94
$ws.Range("U19").Value = 1.91
$ws.Range("V19").Value = 1.72
$ws.Range("W19").Value = 1.33
$ws.Range("X19").Value = 11
$ws.Range("Y19").Value = 8.2
$ws.Range("Z19").Value = 13.5
$ws.Range("AA19").Value = 32
$ws.Range("AB19").Value = 12
$ws.Range("AC19").Value = 7.6
$ws.Range("AD19").Value = 11.5
$ws.Range("AE19").Value = 28
$ws.Range("AF19").Value = 26
$ws.Range("AG19").Value = 16.5
$ws.Range("AH19").Value = 21
$ws.Range("AI19").Value = 70
$ws.Range("AJ19").Value = 95
$ws.Range("AK19").Value = 65
$ws.Range("AL19").Value = 85
$ws.Range("AM19").Value = 170
$ws.Range("AN19").Value = 80
$ws.Range("AO19").Value = 26
